# Apply the "Add files via upload" edit:
#  - 2035_Fx scenario label is repurposed to a new "2035 E" scenario (SB375 ouptut for
#    2035_E_minus_test2), keeping the legacy 2035_F numbers on their own row.
#  - Fill in the previously-blank 2020 and 2025nb (M/N) scenario columns, plus the
#    newly-added 2035 E scenario (Q column) in both data tables.
#  - Extend the lower pivot table (rows 23-29) with the matching 2020 / 2025nb / 2035 E
#    rows, and move the legacy 2035_F row down to row 29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Top table (rows 3-9): Scenario ID / RunEx / VMT / StartEx / Vehicle trips
# ---------------------------------------------------------------------------

# Row 3 - Scenario ID
$ws.Range("M3").Value = 82
$ws.Range("N3").Value = 81
$ws.Range("Q3").Value = 83

# Row 4 - SB375 CO2 RunEx
$ws.Range("M4").Value = 36382.736929999999
$ws.Range("N4").Value = 37296.81624
$ws.Range("Q4").Value = 37398.656410000003

# Row 5 - SB375 VMT
$ws.Range("M5").Value = 77822865.239999995
$ws.Range("N5").Value = 80377022.120000005
$ws.Range("Q5").Value = 81872464.950000003

# Row 7 - SB375 CO2 StartEx
$ws.Range("M7").Value = 1321.2375360000001
$ws.Range("N7").Value = 1424.998083
$ws.Range("Q7").Value = 1577.3501940000001

# Row 8 - SB375 Vehicle Trips
$ws.Range("M8").Value = 14006178.640000001
$ws.Range("N8").Value = 15245816.220000001
$ws.Range("Q8").Value = 16966619.170000002

# ---------------------------------------------------------------------------
# Row 12-18 (second copy of the same table, with text scenario labels)
# ---------------------------------------------------------------------------

# Row 12 - Scenario labels.
# R12 used to be the only cell carrying the "2035_Fx" label; that text is now
# reused for the brand-new "2035 E" scenario in column Q, while R12 goes back
# to the plain "2035_F" label.
$ws.Range("R12").Value = "2035_F"
$ws.Range("M12").Value = 2020
$ws.Range("N12").Value = "2025nb"
$ws.Range("Q12").Value = "2035 E"

# Row 13 - Scenario_ID
$ws.Range("M13").Value = 82
$ws.Range("N13").Value = 81
$ws.Range("Q13").Value = 83

# Row 14 - SB375_VMT
$ws.Range("M14").Value = 77822865.239999995
$ws.Range("N14").Value = 80377022.120000005
$ws.Range("Q14").Value = 81872464.950000003

# Row 15 - SB375_CO2
$ws.Range("M15").Value = 37703.974459999998
$ws.Range("N15").Value = 38721.814319999998
$ws.Range("Q15").Value = 38976.006600000001

# Row 16 - SB375_CO2_runex
$ws.Range("M16").Value = 36382.736929999999
$ws.Range("N16").Value = 37296.81624
$ws.Range("Q16").Value = 37398.656410000003

# Row 17 - SB375_CO2_startEx
$ws.Range("M17").Value = 1321.2375360000001
$ws.Range("N17").Value = 1424.998083
$ws.Range("Q17").Value = 1577.3501940000001

# Row 18 - SB375_vehicle_trip
$ws.Range("M18").Value = 14006178.640000001
$ws.Range("N18").Value = 15245816.220000001
$ws.Range("Q18").Value = 16966619.170000002

# Recalculate so the dependent ratio formulas (rows 6 and 9) pick up real
# values instead of the old #DIV/0! errors.
$excel.Calculate()

# ---------------------------------------------------------------------------
# Lower pivot-style table (rows 22-29)
# ---------------------------------------------------------------------------

# Give the new year rows (23/24) the same left-aligned style as each other.
$ws.Range("K23:K24").HorizontalAlignment = -4131

# Row 24 - new "2020" row
$ws.Range("K24").Value = 2020
$ws.Range("L24").Value = 82
$ws.Range("M24").Value = 77822865.239999995
$ws.Range("N24").Value = 37703.974459999998
$ws.Range("O24").Value = 36382.736929999999
$ws.Range("P24").Value = 1321.2375360000001
$ws.Range("Q24").Value = 14006178.640000001

# Row 25 - new "2025nb" row
$ws.Range("K25").Value = "2025nb"
$ws.Range("L25").Value = 81
$ws.Range("M25").Value = 80377022.120000005
$ws.Range("N25").Value = 38721.814319999998
$ws.Range("O25").Value = 37296.81624
$ws.Range("P25").Value = 1424.998083
$ws.Range("Q25").Value = 15245816.220000001

# Apply the "Comma" style from row 14 onto the rows whose number formatting
# changed (23, 26, 27) or that are brand new (24, 25, 28, 29).
$ws.Range("M14").Copy()
$ws.Range("M23:Q27").PasteSpecial(-4122)
$ws.Range("M28:Q29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 28 used to hold the legacy "2035_Fx" scenario (id 78); it's repointed to
# the new "2035 E" scenario (id 83) with its own, more precise figures.
$ws.Range("K28").Value = "2035_E"
$ws.Range("L28").Value = 83
$ws.Range("M28").Value = 81872464.947604507
$ws.Range("N28").Value = 38976.006602709698
$ws.Range("O28").Value = 37398.6564084731
$ws.Range("P28").Value = 1577.3501942366099
$ws.Range("Q28").Value = 16966619.170381401

# Row 29 - new row holding the legacy "2035_F" scenario numbers that used to
# live on row 28.
$ws.Range("K29").Value = "2035_F"
$ws.Range("L29").Value = 78
$ws.Range("M29").Value = 86586762.788520694
$ws.Range("N29").Value = 41721.8132748598
$ws.Range("O29").Value = 40053.637689939998
$ws.Range("P29").Value = 1668.17558491976
$ws.Range("Q29").Value = 17943573.952111099

$excel.Calculate()

# ---------------------------------------------------------------------------
# Misc: restore the selected cell as last left by the editor.
# ---------------------------------------------------------------------------
$ws.Range("M21").Select()
